$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Rows("12:13").Insert()
$ws.Range("A13").Value = "Output directory:"
